# Fruta / hortaliza, semanal
# Insert two new weekly observations (rows 54-55) into the Mango price
# series, pushing the existing rows 54-95 down to 56-97.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 54:55 - everything from the old row 54 onward
# shifts down by two rows (old 54 -> new 56, ... old 95 -> new 97), and the
# row formatting (e.g. the date style on column D) is carried along by the
# Insert() call automatically.
$ws.Rows("54:55").Insert()

# Populate the two new rows with the new weekly price observations. The
# "constant" columns (market/region/product taxonomy) match every other
# "Sin especificar" / Peru / bandeja row in this block.
$ws.Range("A54").Value2 = 1
$ws.Range("B54").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C54").Value2 = "Arica y Parinacota"
$ws.Range("D54").Value2 = 44484
$ws.Range("E54").Value2 = 15
$ws.Range("F54").Value2 = "Fruta"
$ws.Range("G54").Value2 = 100108
$ws.Range("H54").Value2 = "Tropicales y subtropicales"
$ws.Range("I54").Value2 = 100108002
$ws.Range("J54").Value2 = "Mango"
$ws.Range("K54").Value2 = "Sin especificar"
$ws.Range("L54").Value2 = "Especial"
$ws.Range("M54").Value2 = 456
$ws.Range("N54").Value2 = 4500
$ws.Range("O54").Value2 = 5000
$ws.Range("P54").Value2 = 4750
$ws.Range("Q54").Value2 = "$/bandeja 4 kilos"
$ws.Range("R54").Value2 = "Perú"
$ws.Range("S54").Value2 = 1188
$ws.Range("T54").Value2 = 4

$ws.Range("A55").Value2 = 1
$ws.Range("B55").Value2 = "Agrícola del Norte S.A. de Arica"
$ws.Range("C55").Value2 = "Arica y Parinacota"
$ws.Range("D55").Value2 = 44484
$ws.Range("E55").Value2 = 15
$ws.Range("F55").Value2 = "Fruta"
$ws.Range("G55").Value2 = 100108
$ws.Range("H55").Value2 = "Tropicales y subtropicales"
$ws.Range("I55").Value2 = 100108002
$ws.Range("J55").Value2 = "Mango"
$ws.Range("K55").Value2 = "Sin especificar"
$ws.Range("L55").Value2 = "Primera"
$ws.Range("M55").Value2 = 456
$ws.Range("N55").Value2 = 4500
$ws.Range("O55").Value2 = 5000
$ws.Range("P55").Value2 = 4750
$ws.Range("Q55").Value2 = "$/bandeja 4 kilos"
$ws.Range("R55").Value2 = "Perú"
$ws.Range("S55").Value2 = 1188
$ws.Range("T55").Value2 = 4
